# Reproduce the authored change:
#  - B1 gets a new value "Test1" (second shared string, dimension A1 -> A1:B1)
#  - Active selection moves from A2 to B2
#  - Workbook gains two custom (personal) views, one per collaborator,
#    as created when the file was turned into a shared workbook
#    ("APS_DEV_01 - Affichage personnalisé" and
#    "Alexandre Baranger - Affichage personnalisé")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- data edit: B1 = "Test1" ---
$ws.Range("B1").Value = "Test1"

# --- selection moves to B2 ---
$ws.Range("B2").Select() | Out-Null

# --- custom workbook views added (order matches the authored file) ---
$wb.CustomViews.Add("APS_DEV_01 - Affichage personnalisé") | Out-Null
$wb.CustomViews.Add("Alexandre Baranger - Affichage personnalisé") | Out-Null
